# Se agregan las tareas del dia - se arregla consulta de producto por codigo
#
# Adds six new task rows (rows 52-57) to the bottom of the task list on
# Hoja1, moves the existing "Ivan: preguntar reportes..." row from row 55
# down to row 59 (leaving row 58 blank), and updates the sheet selection
# to reflect where the user ended up after the edits (cell B56).
#
# NOTE on shared-string insertion order: Excel (and this COM layer) assigns
# new sharedStrings.xml entries in first-write order, not sheet/row order.
# To reproduce the exact target shared-string indices the new A-column
# labels are written in the same order the original author must have typed
# them (A52, A53, A55, A54, A56, A57) before the B-column responsible
# names are filled in.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New task descriptions (column A) - written in "first seen" order so the
# shared string table comes out in the same order as the target workbook.
$ws.Cells.Item(52, 1).Value = "Estetica - botones imágenes"
$ws.Cells.Item(53, 1).Value = "Estetica - etiquetas - tooltips - mensajes - mensajes de errores - etc"
$ws.Cells.Item(55, 1).Value = "Logueo de aplicación (configurable)"
$ws.Cells.Item(54, 1).Value = "Reportes - ruta - estitca"
$ws.Cells.Item(56, 1).Value = "Reunion con Ivan y Josefina - consutlas"
$ws.Cells.Item(57, 1).Value = "Impresora - carga de datos y factura"

# Responsible person (column B) for the new rows.
$ws.Cells.Item(54, 2).Value = "Agustina"
$ws.Cells.Item(55, 2).Value = "Lucas"
$ws.Cells.Item(56, 2).Value = "Lucas/Agustina"
$ws.Cells.Item(57, 2).Value = "Lucas/Agustina"

# Row 55 used to hold "Ivan: preguntar reportes..." - that task gets
# bumped down to row 59 (row 58 stays empty) to make room for the new
# "Logueo de aplicación" row above.
$ws.Cells.Item(59, 1).Value = "Ivan: preguntar reportes - preguntar autorizacion requerida en que funciones - preguntar login"

# Reflect the final cursor/selection position left by the edit.
$ws.Activate()
$ws.Range("B56").Select()
